$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.758.68'
$ws.Range('E2').Value = '  -10.19%  '
$ws.Range('D3').Value = '2.390.61'
$ws.Range('E3').Value = '  -16.57%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '461.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -10.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.44'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -7.05%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.993'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.482'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -9.17%  '
$ws.Range('D9').Value = '2.403.37'
$ws.Range('E9').Value = '  -15.97%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0939'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.32'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -12.10%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.316'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -9.75%  '
$ws.Range('E13').Value = '  -4.38%  '
$ws.Range('D14').Value = '2.769.93'
$ws.Range('E14').Value = '  -17.60%  '
$ws.Range('D15').Value = '53.696.56'
$ws.Range('E15').Value = '  -10.71%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.45'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -11.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000125'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -8.29%  '
$ws.Range('D18').Value = '2.379.23'
$ws.Range('E18').Value = '  -17.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.18'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -12.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '308.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -11.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.22'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -17.89%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.02'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.67'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -15.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '55.43'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -12.91%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.380'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -12.87%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.154'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -12.01%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.443.12'
$ws.Range('E29').Value = '  -18.77%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.05'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.994'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '0.0₃0701'
$ws.Range('E32').Value = '  -16.15%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '147.10'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.56'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -8.19%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.38'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -15.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.96'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -8.63%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.43'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -19.40%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.03'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -11.14%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '33.28'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -10.20%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.992'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.768'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -20.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.599'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.27'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -8.65%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0522'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -7.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.12'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('D46').Value = '1.950.87'
$ws.Range('E46').Value = '  -12.98%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.21'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -14.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0213'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -7.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0852'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.26'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -9.58%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.40'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -17.32%  '
